$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "6. Ideja je uklonjena i korisnik je vraćen na naslovnu stranu."
# becomes
# "6. Ideja je uklonjena, kao i svako predviđanje koje je nastalo kao odgovor
#  na ideju, i korisnik je vraćen na naslovnu stranu."
$d.Content.Find.Execute(
    "6. Ideja je uklonjena i korisnik je vraćen na naslovnu stranu.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "6. Ideja je uklonjena, kao i svako predviđanje koje je nastalo kao odgovor na ideju, i korisnik je vraćen na naslovnu stranu.",
    2)

# Give the newly inserted "predviđanje ... ideju," fragment the Serbian
# (Latin) language tag, matching the author's edit.
$r1 = $d.Content
$r1.Find.Execute("đanje koje je nastalo kao odgovor na ideju,", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.LanguageID = "sr-Latn-RS"

# --- Change 2 -------------------------------------------------------------
# Append a sentence to the end of the "Posledice" paragraph.
$r2 = $d.Content
$r2.Find.Execute(
    "Posledica brisanja ideje je njeno uklanjanje iz baze podataka, kao i automatsko uklanjanje iz sistema svih predviđanja nastalih kao odgovor na datu ideju.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertAfter(" Takođe je, eventualno, ako je predviđanje završeno, potrebno ažurirati skor svim korisnicima koji su odgovorili na predviđanja.")
